$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.266377860180107
$ws.Range("D2").Value = 0.1022916732020462
$ws.Range("E2").Value = 0.08516941909815312
$ws.Range("F2").Value = 2.436614961427182
$ws.Range("G2").Value = 2.533292934628946
$ws.Range("H2").Value = 1.863095474905975
$ws.Range("L2").Value = 0.07448513552507308
$ws.Range("M2").Value = 2.122762403723868
$ws.Range("N2").Value = 1.541936691234667
$ws.Range("C3").Value = 0.2689043029714995
$ws.Range("D3").Value = 0.1014319877055812
$ws.Range("E3").Value = 0.08653199089029151
$ws.Range("F3").Value = 2.347201262919498
$ws.Range("G3").Value = 2.425336855193791
$ws.Range("H3").Value = 1.822689697272722
$ws.Range("L3").Value = 0.07322391023355124
$ws.Range("M3").Value = 1.940322024169575
$ws.Range("N3").Value = 1.44143269200066
$ws.Range("C4").Value = 0.2706335220392333
$ws.Range("D4").Value = 0.1010012580974688
$ws.Range("E4").Value = 0.08741164768381937
$ws.Range("F4").Value = 2.294427139550351
$ws.Range("G4").Value = 2.361357575447983
$ws.Range("H4").Value = 1.799392259840204
$ws.Range("L4").Value = 0.07245961574455606
$ws.Range("M4").Value = 1.828709854513932
$ws.Range("N4").Value = 1.380082847838764
$ws.Range("C5").Value = 0.2713827151744539
$ws.Range("D5").Value = 0.1008498367917383
$ws.Range("E5").Value = 0.08778093051892899
$ws.Range("F5").Value = 2.273447776089057
$ws.Range("G5").Value = 2.335855500065634
$ws.Range("H5").Value = 1.790274081009727
$ws.Range("L5").Value = 0.07215076241349649
$ws.Range("M5").Value = 1.783328225325292
$ws.Range("N5").Value = 1.355175481027828
$ws.Range("C6").Value = 0.2715098005036793
$ws.Range("D6").Value = 0.1008261400779276
$ws.Range("E6").Value = 0.0878429027366896
$ws.Range("F6").Value = 2.269995733970177
$ws.Range("G6").Value = 2.331655039992313
$ws.Range("H6").Value = 1.788782577971176
$ws.Range("L6").Value = 0.07209963659519403
$ws.Range("M6").Value = 1.775798725071979
$ws.Range("N6").Value = 1.351045339948399
$ws.Range("C7").Value = 0.2706434459514213
$ws.Range("D7").Value = 0.100999118793375
$ws.Range("E7").Value = 0.08741658417601661
$ws.Range("F7").Value = 2.294142083305019
$ws.Range("G7").Value = 2.361011350595163
$ws.Range("H7").Value = 1.79926777345301
$ws.Range("L7").Value = 0.07245543981914793
$ws.Range("M7").Value = 1.828097413460313
$ws.Range("N7").Value = 1.379746557464188
$ws.Range("C8").Value = 0.2672118994173687
$ws.Range("D8").Value = 0.1019749030892427
$ws.Range("E8").Value = 0.08563029946008172
$ws.Range("F8").Value = 2.40533884068526
$ws.Range("G8").Value = 2.495584977955531
$ws.Range("H8").Value = 1.848847049858676
$ws.Range("L8").Value = 0.07404820515161248
$ws.Range("M8").Value = 2.059771745073959
$ws.Range("N8").Value = 1.507209734598149
$ws.Range("C9").Value = 0.2619053908150732
$ws.Range("D9").Value = 0.1046737890596035
$ws.Range("E9").Value = 0.08246902750140439
$ws.Range("F9").Value = 2.64065259201567
$ws.Range("G9").Value = 2.778251224005544
$ws.Range("H9").Value = 1.958279468790806
$ws.Range("L9").Value = 0.07724919980169886
$ws.Range("M9").Value = 2.517392233225877
$ws.Range("N9").Value = 1.75992046572614
$ws.Range("C10").Value = 0.258889057753521
$ws.Range("D10").Value = 0.1071565716895293
$ws.Range("E10").Value = 0.08035489119832295
$ws.Range("F10").Value = 2.8246213403539
$ws.Range("G10").Value = 2.998049845067328
$ws.Range("H10").Value = 2.046430059291481
$ws.Range("L10").Value = 0.07964505071615946
$ws.Range("M10").Value = 2.855774278434097
$ws.Range("N10").Value = 1.947155578148198
$ws.Range("C11").Value = 0.2577117996959402
$ws.Range("D11").Value = 0.1083993215766696
$ws.Range("E11").Value = 0.079438422332613
$ws.Range("F11").Value = 2.910846403123486
$ws.Range("G11").Value = 3.100827179895248
$ws.Range("H11").Value = 2.088284246704404
$ws.Range("L11").Value = 0.08074388846999625
$ws.Range("M11").Value = 3.010221289522207
$ws.Range("N11").Value = 2.032653215215078
$ws.Range("C12").Value = 0.2572943117838946
$ws.Range("D12").Value = 0.1088866046441268
$ws.Range("E12").Value = 0.07909789518356991
$ws.Range("F12").Value = 2.943872559896988
$ws.Range("G12").Value = 3.140159796490479
$ws.Range("H12").Value = 2.104391178681965
$ws.Range("L12").Value = 0.08116121611568161
$ws.Range("M12").Value = 3.068782935573637
$ws.Range("N12").Value = 2.065073364915975
$ws.Range("C13").Value = 0.2573829613079681
$ws.Range("D13").Value = 0.1087809114022917
$ws.Range("E13").Value = 0.07917094365062116
$ws.Range("F13").Value = 2.936742966974833
$ws.Range("G13").Value = 3.131670249888714
$ws.Range("H13").Value = 2.100910715241923
$ws.Range("L13").Value = 0.08107128385472606
$ws.Range("M13").Value = 3.056167230445595
$ws.Range("N13").Value = 2.058089184882135
$ws.Range("C14").Value = 0.2576768835718042
$ws.Range("D14").Value = 0.1084390741009287
$ws.Range("E14").Value = 0.07941027622959496
$ws.Range("F14").Value = 2.913555920372772
$ws.Range("G14").Value = 3.104054750238333
$ws.Range("H14").Value = 2.089604178220327
$ws.Range("L14").Value = 0.08077819815972731
$ws.Range("M14").Value = 3.015037659005969
$ws.Range("N14").Value = 2.035319571454863
$ws.Range("C15").Value = 0.2578606155031622
$ws.Range("D15").Value = 0.108231872797063
$ws.Range("E15").Value = 0.07955772373949177
$ws.Range("F15").Value = 2.899402263737926
$ws.Range("G15").Value = 3.087193611173802
$ws.Range("H15").Value = 2.082712316208926
$ws.Range("L15").Value = 0.08059883185010364
$ws.Range("M15").Value = 2.989854551702678
$ws.Range("N15").Value = 2.02137819058089
$ws.Range("C16").Value = 0.2589699405614141
$ws.Range("D16").Value = 0.1070776680771246
$ws.Range("E16").Value = 0.0804156963233007
$ws.Range("F16").Value = 2.819038209281985
$ws.Range("G16").Value = 2.991390253606255
$ws.Range("H16").Value = 2.043730546610561
$ws.Range("L16").Value = 0.07957341361718306
$ws.Range("M16").Value = 2.84569132402541
$ws.Range("N16").Value = 1.941574403025243
$ws.Range("C17").Value = 0.2597005946676632
$ws.Range("D17").Value = 0.1063988960452775
$ws.Range("E17").Value = 0.0809536352155602
$ws.Range("F17").Value = 2.770394097286754
$ws.Range("G17").Value = 2.93334091623052
$ws.Range("H17").Value = 2.02026956656249
$ws.Range("L17").Value = 0.0789465990071605
$ws.Range("M17").Value = 2.757385080424513
$ws.Range("N17").Value = 1.892698432424453
$ws.Range("C18").Value = 0.2601391735297085
$ws.Range("D18").Value = 0.1060191268947079
$ws.Range("E18").Value = 0.08126730306006635
$ws.Range("F18").Value = 2.742653461000884
$ws.Range("G18").Value = 2.900214337336593
$ws.Range("H18").Value = 2.006940197127221
$ws.Range("L18").Value = 0.07858691846512755
$ws.Range("M18").Value = 2.706642051392919
$ws.Range("N18").Value = 1.864616853602627
$ws.Range("C19").Value = 0.260290807581697
$ws.Range("D19").Value = 0.1058923598644839
$ws.Range("E19").Value = 0.08137423686552392
$ws.Range("F19").Value = 2.733301568982029
$ws.Range("G19").Value = 2.889042864562782
$ws.Range("H19").Value = 2.002455237054903
$ws.Range("L19").Value = 0.07846528423174703
$ws.Range("M19").Value = 2.689469591405526
$ws.Range("N19").Value = 1.855114242984854
$ws.Range("C20").Value = 0.2596209165493377
$ws.Range("D20").Value = 0.1064700482595242
$ws.Range("E20").Value = 0.08089592976896487
$ws.Range("F20").Value = 2.775547625994989
$ws.Range("G20").Value = 2.939493180039051
$ws.Range("H20").Value = 2.022749936412822
$ws.Range("L20").Value = 0.07901323734508026
$ws.Range("M20").Value = 2.766780410329915
$ws.Range("N20").Value = 1.897898215375875
$ws.Range("C21").Value = 0.2575897806953265
$ws.Range("D21").Value = 0.1085390241764941
$ws.Range("E21").Value = 0.07933980141602248
$ws.Range("F21").Value = 2.920356272162849
$ws.Range("G21").Value = 3.112154785070857
$ws.Range("H21").Value = 2.092918145443775
$ws.Range("L21").Value = 0.08086425190214896
$ws.Range("M21").Value = 3.027116335612163
$ws.Range("N21").Value = 2.042006382610339
$ws.Range("C22").Value = 0.2564274921204657
$ws.Range("D22").Value = 0.1099886328152309
$ws.Range("E22").Value = 0.07836078892418463
$ws.Range("F22").Value = 3.01718525756101
$ws.Range("G22").Value = 3.22741254109053
$ws.Range("H22").Value = 2.140281582570822
$ws.Range("L22").Value = 0.08208108705827044
$ws.Range("M22").Value = 3.197704354941635
$ws.Range("N22").Value = 2.136445066950387
$ws.Range("C23").Value = 0.2570326162253593
$ws.Range("D23").Value = 0.1092059062675474
$ws.Range("E23").Value = 0.07887982488053424
$ws.Range("F23").Value = 2.965302283074237
$ws.Range("G23").Value = 3.16567247508118
$ws.Range("H23").Value = 2.114863336579447
$ws.Range("L23").Value = 0.08143101221343585
$ws.Range("M23").Value = 3.106617104178724
$ws.Range("N23").Value = 2.086018764779794
$ws.Range("C24").Value = 0.259656881374994
$ws.Range("D24").Value = 0.1064378477739183
$ws.Range("E24").Value = 0.08092200469848798
$ws.Range("F24").Value = 2.773217013987107
$ws.Range("G24").Value = 2.93671097387022
$ws.Range("H24").Value = 2.021628067420522
$ws.Range("L24").Value = 0.07898310801856923
$ws.Range("M24").Value = 2.762532702677674
$ws.Range("N24").Value = 1.895547337661725
$ws.Range("C25").Value = 0.2631870981078492
$ws.Range("D25").Value = 0.1038574353849953
$ws.Range("E25").Value = 0.08328764407238998
$ws.Range("F25").Value = 2.575089023843844
$ws.Range("G25").Value = 2.699703667517099
$ws.Range("H25").Value = 1.927336742911564
$ws.Range("L25").Value = 0.076375272817927
$ws.Range("M25").Value = 2.393224463199772
$ws.Range("N25").Value = 1.691272285562036
